# Minor changes to the word file.
#
# - Remove the two blank paragraphs immediately before the
#   "Then type new template blank..." paragraph.
# - Remove the two blank paragraphs immediately after it.
# - Move the (hidden) "_GoBack" bookmark from the end of that paragraph's
#   text to the beginning of it.

$d = $word.ActiveDocument

# Locate the target paragraph by its distinctive text (robust to index
# shifts caused by earlier edits elsewhere in the body). We re-look-up by
# index (rather than keep a Paragraph object reference) because paragraph
# object references get reseated by position once the document is
# mutated.
$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Then type new template blank*") {
        $targetIdx = $i
        break
    }
}

# Remove the two blank paragraphs that follow the target paragraph first
# -- this does not change the target paragraph's own index.
$d.Paragraphs.Item($targetIdx + 1).Range.Delete()
$d.Paragraphs.Item($targetIdx + 1).Range.Delete()

# Remove the two blank paragraphs that precede the target paragraph.
# Each deletion shifts the target paragraph's index down by one.
$d.Paragraphs.Item($targetIdx - 1).Range.Delete()
$targetIdx = $targetIdx - 1
$d.Paragraphs.Item($targetIdx - 1).Range.Delete()
$targetIdx = $targetIdx - 1

# Move the "_GoBack" bookmark from the end of the paragraph text to the
# start of it.
$paraStart = $d.Paragraphs.Item($targetIdx).Range.Start
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()
$newRange = $d.Range($paraStart, $paraStart)
$d.Bookmarks.Add("_GoBack", $newRange)
